# Inserts a new data row at row 239 (pushing the existing rows 239-334 down
# to 240-335) and populates the new row with the new "Feria Lagunitas de
# Puerto Montt - Betarraga" price observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 239..334 down to 240..335, leaving a blank row at 239.
$ws.Rows("239").Insert()

# Populate the newly inserted row 239 with the new record's data.
$ws.Range("A239").Value = 4
$ws.Range("B239").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C239").Value = "Los Lagos"
$ws.Range("D239").Value = 44784
$ws.Range("E239").Value = 10
$ws.Range("F239").Value = 100114014
$ws.Range("G239").Value = "Betarraga"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 500
$ws.Range("K239").Value = 1000
$ws.Range("L239").Value = 1200
$ws.Range("M239").Value = 1100
$ws.Range("N239").Value = "$/paquete 5 unidades"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 220
$ws.Range("Q239").Value = 5
$ws.Range("R239").Value = "Hortaliza"
